# Add the "Project" value ("100 Home Project") to column C for the rows
# that were missing it on the "Location Adjustment" sheet (rows 2, 4, 6, 8, 10),
# matching rows 3, 5, 7, 9, 11 which already carry that value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Location Adjustment")

$rows = @(2, 4, 6, 8, 10)
foreach ($r in $rows) {
    $ws.Range("C$r").Value = "100 Home Project"
}
